$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.188222646713257
$ws.Range("B1").Value = 1.307642340660095
$ws.Range("C1").Value = 1.553058981895447
$ws.Range("D1").Value = 2.722550630569458
$ws.Range("E1").Value = 15
